$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column A to fit the longer task names ---
# ColumnWidth (character units) differs from the stored XML "width" by a
# fixed offset on this engine; subtract it so the saved width comes out to
# exactly 44.
$ws.Columns.Item(1).ColumnWidth = 43.166666666666664

# --- Task list: reordered / renamed rows with new Start Date & Days to
#     Complete columns (B & C) for the initial time estimations ---
$tasks = @(
    @{ Row = 2;  Name = "Integrate old Lab Work (GJ)";            Start = 43410; Days = 2 },
    @{ Row = 3;  Name = "Client GUI";                              Start = 43412; Days = 2 },
    @{ Row = 4;  Name = "Server GUI";                               Start = 43412; Days = 2 },
    @{ Row = 5;  Name = "Start Game GUI and functionality (AJ)";   Start = 43412; Days = 2 },
    @{ Row = 6;  Name = "Gameplay GUI";                             Start = 43412; Days = 4 },
    @{ Row = 7;  Name = "Bet GUI";                                  Start = 43412; Days = 4 },
    @{ Row = 8;  Name = "Database Creation";                        Start = 43412; Days = 2 },
    @{ Row = 9;  Name = "Database Communication";                   Start = 43414; Days = 4 },
    @{ Row = 10; Name = "Gameplay Functionality";                   Start = 43416; Days = 9 },
    @{ Row = 11; Name = "Bet Functionality";                        Start = 43416; Days = 8 },
    @{ Row = 12; Name = "Unit Testing (all)";                       Start = 43425; Days = 3 },
    @{ Row = 13; Name = "Integration Testing (all)";                Start = 43428; Days = 5 },
    @{ Row = 14; Name = "Final Testing (all)";                      Start = 43433; Days = 4 }
)

foreach ($t in $tasks) {
    $r = $t.Row
    $ws.Cells.Item($r, 1).Value = $t.Name
    $ws.Cells.Item($r, 2).Value = $t.Start
    $ws.Cells.Item($r, 2).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 3).Value = $t.Days
}

# --- Leave behind the sort remnant (Data > Sort by Start Date) ---
$ws.Sort.SortFields.Add($ws.Range("B2:B14"))
$ws.Sort.SetRange($ws.Range("A2:C13"))
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.Apply()

# --- Page orientation set to portrait via Page Setup ---
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# --- Final selection rests on the empty row below the table ---
$ws.Range("A15").Select() | Out-Null
